$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (rich-text cells): edit in place via Characters to preserve surrounding text ---
# A8 = "Volume 29   Number  38" -> "...Number  40"  (chars 21-22 are "38")
$ws.Range("A8").Characters(21, 2).Text = "40"

# C9 = "Report Covering the Week  9/19/2022  Through  9/25/2022"
#      -> "...10/3/2022  Through  10/9/2022" (chars 27-35 and 47-55 are the dates)
$ws.Range("C9").Characters(27, 9).Text = "10/3/2022"
$ws.Range("C9").Characters(47, 9).Text = "10/9/2022"


# --- Style-changing cells: copy style (and text value where applicable) from stable reference cells ---
$ws.Range("D22").Copy($ws.Range("D16"))
$ws.Range("E22").Copy($ws.Range("E16"))
$ws.Range("D22").Copy($ws.Range("C22"))
$ws.Range("D22").Copy($ws.Range("D30"))
$ws.Range("E22").Copy($ws.Range("E30"))
$ws.Range("I14").Copy($ws.Range("C15"))
$ws.Range("C15").Value = 1
$ws.Range("I14").Copy($ws.Range("F15"))
$ws.Range("F15").Value = 1
$ws.Range("I14").Copy($ws.Range("C26"))
$ws.Range("C26").Value = 2
$ws.Range("I14").Copy($ws.Range("F26"))
$ws.Range("F26").Value = 2

# --- Plain value updates (style unchanged) ---
$ws.Range("G15").Value = 1
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 11
$ws.Range("K15").Value = 22.222222222222
$ws.Range("L15").Value = 120
$ws.Range("M15").Value = 10
$ws.Range("N15").Value = 0
$ws.Range("C16").Value = 6
$ws.Range("F16").Value = 28
$ws.Range("H16").Value = 55.555555555555
$ws.Range("I16").Value = 200
$ws.Range("J16").Value = 135
$ws.Range("K16").Value = 48.148148148148
$ws.Range("L16").Value = 92.307692307692
$ws.Range("M16").Value = 86.915887850467
$ws.Range("N16").Value = -70.104633781763
$ws.Range("C17").Value = 2
$ws.Range("D17").Value = 1
$ws.Range("E17").Value = 100
$ws.Range("F17").Value = 15
$ws.Range("G17").Value = 17
$ws.Range("H17").Value = -11.764705882352
$ws.Range("I17").Value = 148
$ws.Range("J17").Value = 140
$ws.Range("K17").Value = 5.714285714285
$ws.Range("L17").Value = 21.311475409836
$ws.Range("M17").Value = 108.450704225352
$ws.Range("N17").Value = -38.333333333333
$ws.Range("C18").Value = 6
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = 100
$ws.Range("F18").Value = 29
$ws.Range("G18").Value = 30
$ws.Range("H18").Value = -3.333333333333
$ws.Range("I18").Value = 341
$ws.Range("J18").Value = 176
$ws.Range("K18").Value = 93.75
$ws.Range("L18").Value = 28.679245283018
$ws.Range("M18").Value = 143.571428571429
$ws.Range("N18").Value = -47.699386503067
$ws.Range("C19").Value = 22
$ws.Range("D19").Value = 21
$ws.Range("E19").Value = 4.761904761904
$ws.Range("F19").Value = 120
$ws.Range("G19").Value = 99
$ws.Range("H19").Value = 21.212121212121
$ws.Range("I19").Value = 1029
$ws.Range("J19").Value = 561
$ws.Range("K19").Value = 83.422459893048
$ws.Range("L19").Value = 93.785310734463
$ws.Range("M19").Value = 27.509293680297
$ws.Range("N19").Value = -47.122302158273
$ws.Range("F20").Value = 4
$ws.Range("H20").Value = -20
$ws.Range("I20").Value = 50
$ws.Range("J20").Value = 40
$ws.Range("K20").Value = 25
$ws.Range("L20").Value = 19.047619047619
$ws.Range("M20").Value = 51.515151515151
$ws.Range("N20").Value = -91.304347826087
$ws.Range("C21").Value = 38
$ws.Range("D21").Value = 26
$ws.Range("E21").Value = 46.153846153846
$ws.Range("F21").Value = 197
$ws.Range("G21").Value = 170
$ws.Range("H21").Value = 15.882352941176
$ws.Range("I21").Value = 1781
$ws.Range("J21").Value = 1061
$ws.Range("K21").Value = 67.860508953817
$ws.Range("L21").Value = 66.448598130841
$ws.Range("M21").Value = 52.482876712328
$ws.Range("N21").Value = -56.497313141182
$ws.Range("F22").Value = 4
$ws.Range("G22").Value = 7
$ws.Range("H22").Value = -42.857142857142
$ws.Range("J22").Value = 26
$ws.Range("K22").Value = 26.923076923076
$ws.Range("L22").Value = 32
$ws.Range("M22").Value = -15.384615384615
$ws.Range("C24").Value = 30
$ws.Range("D24").Value = 30
$ws.Range("E24").Value = 0
$ws.Range("F24").Value = 151
$ws.Range("G24").Value = 118
$ws.Range("H24").Value = 27.966101694915
$ws.Range("I24").Value = 1526
$ws.Range("J24").Value = 979
$ws.Range("K24").Value = 55.873340143003
$ws.Range("L24").Value = 53.83064516129
$ws.Range("M24").Value = 30.987124463519
$ws.Range("C25").Value = 11
$ws.Range("D25").Value = 3
$ws.Range("E25").Value = 266.666666666667
$ws.Range("F25").Value = 38
$ws.Range("G25").Value = 29
$ws.Range("H25").Value = 31.03448275862
$ws.Range("I25").Value = 317
$ws.Range("J25").Value = 237
$ws.Range("K25").Value = 33.755274261603
$ws.Range("L25").Value = 58.5
$ws.Range("M25").Value = 64.248704663212
$ws.Range("G26").Value = 1
$ws.Range("H26").Value = 100
$ws.Range("I26").Value = 16
$ws.Range("K26").Value = 45.454545454545
$ws.Range("L26").Value = 128.571428571429
$ws.Range("C27").Value = 2
$ws.Range("D27").Value = 2
$ws.Range("E27").Value = 0
$ws.Range("G27").Value = 11
$ws.Range("H27").Value = -9.090909090909
$ws.Range("I27").Value = 62
$ws.Range("J27").Value = 51
$ws.Range("K27").Value = 21.56862745098
$ws.Range("L27").Value = 93.75
$ws.Range("G30").Value = 2
$ws.Range("L30").Value = 75
